# Auto-generated Excel COM-interop script to add "PO Forecast" sheet
# and rename header columns in existing sheets per the target diff.

$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet at the end ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row ---
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# --- 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$wsForecast.Cells.Item(2,1).Value = 45032.99999999999
$wsForecast.Cells.Item(2,2).Value = 71
$wsForecast.Cells.Item(2,3).Value = -46.43791321538216
$wsForecast.Cells.Item(2,4).Value = 186.8032931752883
$wsForecast.Cells.Item(3,1).Value = 45039.99999999999
$wsForecast.Cells.Item(3,2).Value = 72
$wsForecast.Cells.Item(3,3).Value = -47.19432322034194
$wsForecast.Cells.Item(3,4).Value = 184.7629172564357
$wsForecast.Cells.Item(4,1).Value = 45060.99999999999
$wsForecast.Cells.Item(4,2).Value = 74
$wsForecast.Cells.Item(4,3).Value = -45.18094535741655
$wsForecast.Cells.Item(4,4).Value = 191.2916433627699
$wsForecast.Cells.Item(5,1).Value = 45067.99999999999
$wsForecast.Cells.Item(5,2).Value = 74
$wsForecast.Cells.Item(5,3).Value = -44.26463968047852
$wsForecast.Cells.Item(5,4).Value = 182.4382338150261
$wsForecast.Cells.Item(6,1).Value = 45074.99999999999
$wsForecast.Cells.Item(6,2).Value = 75
$wsForecast.Cells.Item(6,3).Value = -41.87553885725413
$wsForecast.Cells.Item(6,4).Value = 191.9873557892802
$wsForecast.Cells.Item(7,1).Value = 45081.99999999999
$wsForecast.Cells.Item(7,2).Value = 76
$wsForecast.Cells.Item(7,3).Value = -40.92291525188291
$wsForecast.Cells.Item(7,4).Value = 187.4590960783607
$wsForecast.Cells.Item(8,1).Value = 45088.99999999999
$wsForecast.Cells.Item(8,2).Value = 77
$wsForecast.Cells.Item(8,3).Value = -43.11978877853108
$wsForecast.Cells.Item(8,4).Value = 190.3623839241408
$wsForecast.Cells.Item(9,1).Value = 45095.99999999999
$wsForecast.Cells.Item(9,2).Value = 77
$wsForecast.Cells.Item(9,3).Value = -45.99512965469623
$wsForecast.Cells.Item(9,4).Value = 200.9562215662505
$wsForecast.Cells.Item(10,1).Value = 45102.99999999999
$wsForecast.Cells.Item(10,2).Value = 78
$wsForecast.Cells.Item(10,3).Value = -38.64280075469556
$wsForecast.Cells.Item(10,4).Value = 199.7358625888254
$wsForecast.Cells.Item(11,1).Value = 45109.99999999999
$wsForecast.Cells.Item(11,2).Value = 79
$wsForecast.Cells.Item(11,3).Value = -34.45153360992744
$wsForecast.Cells.Item(11,4).Value = 205.6698225259995
$wsForecast.Cells.Item(12,1).Value = 45116.99999999999
$wsForecast.Cells.Item(12,2).Value = 79
$wsForecast.Cells.Item(12,3).Value = -36.73761727714172
$wsForecast.Cells.Item(12,4).Value = 198.8136424741582
$wsForecast.Cells.Item(13,1).Value = 45137.99999999999
$wsForecast.Cells.Item(13,2).Value = 82
$wsForecast.Cells.Item(13,3).Value = -28.36248678109474
$wsForecast.Cells.Item(13,4).Value = 205.5667284349158
$wsForecast.Cells.Item(14,1).Value = 45144.99999999999
$wsForecast.Cells.Item(14,2).Value = 82
$wsForecast.Cells.Item(14,3).Value = -30.24139107174077
$wsForecast.Cells.Item(14,4).Value = 195.1890778648966
$wsForecast.Cells.Item(15,1).Value = 45151.99999999999
$wsForecast.Cells.Item(15,2).Value = 83
$wsForecast.Cells.Item(15,3).Value = -36.42454421428962
$wsForecast.Cells.Item(15,4).Value = 191.2273434570222
$wsForecast.Cells.Item(16,1).Value = 45158.99999999999
$wsForecast.Cells.Item(16,2).Value = 84
$wsForecast.Cells.Item(16,3).Value = -41.08846846078413
$wsForecast.Cells.Item(16,4).Value = 208.7614765301728
$wsForecast.Cells.Item(17,1).Value = 45179.99999999999
$wsForecast.Cells.Item(17,2).Value = 86
$wsForecast.Cells.Item(17,3).Value = -27.12231961905076
$wsForecast.Cells.Item(17,4).Value = 213.0471396506499
$wsForecast.Cells.Item(18,1).Value = 45193.99999999999
$wsForecast.Cells.Item(18,2).Value = 87
$wsForecast.Cells.Item(18,3).Value = -32.53742900003123
$wsForecast.Cells.Item(18,4).Value = 203.7267929257157
$wsForecast.Cells.Item(19,1).Value = 45200.99999999999
$wsForecast.Cells.Item(19,2).Value = 88
$wsForecast.Cells.Item(19,3).Value = -31.64938357536544
$wsForecast.Cells.Item(19,4).Value = 207.1028919607417
$wsForecast.Cells.Item(20,1).Value = 45207.99999999999
$wsForecast.Cells.Item(20,2).Value = 89
$wsForecast.Cells.Item(20,3).Value = -31.08146642435455
$wsForecast.Cells.Item(20,4).Value = 208.6292971286488
$wsForecast.Cells.Item(21,1).Value = 45214.99999999999
$wsForecast.Cells.Item(21,2).Value = 90
$wsForecast.Cells.Item(21,3).Value = -25.76848197306997
$wsForecast.Cells.Item(21,4).Value = 196.872086750384
$wsForecast.Cells.Item(22,1).Value = 45221.99999999999
$wsForecast.Cells.Item(22,2).Value = 90
$wsForecast.Cells.Item(22,3).Value = -29.41677815663937
$wsForecast.Cells.Item(22,4).Value = 202.7856818146017
$wsForecast.Cells.Item(23,1).Value = 45228.99999999999
$wsForecast.Cells.Item(23,2).Value = 91
$wsForecast.Cells.Item(23,3).Value = -18.61903283570707
$wsForecast.Cells.Item(23,4).Value = 213.0233174090731
$wsForecast.Cells.Item(24,1).Value = 45235.99999999999
$wsForecast.Cells.Item(24,2).Value = 92
$wsForecast.Cells.Item(24,3).Value = -30.01725326337007
$wsForecast.Cells.Item(24,4).Value = 202.3371353393307
$wsForecast.Cells.Item(25,1).Value = 45242.99999999999
$wsForecast.Cells.Item(25,2).Value = 92
$wsForecast.Cells.Item(25,3).Value = -29.27895052433436
$wsForecast.Cells.Item(25,4).Value = 207.8234325502558
$wsForecast.Cells.Item(26,1).Value = 45249.99999999999
$wsForecast.Cells.Item(26,2).Value = 93
$wsForecast.Cells.Item(26,3).Value = -31.56327086850019
$wsForecast.Cells.Item(26,4).Value = 203.5413033636013
$wsForecast.Cells.Item(27,1).Value = 45256.99999999999
$wsForecast.Cells.Item(27,2).Value = 94
$wsForecast.Cells.Item(27,3).Value = -32.46626827078989
$wsForecast.Cells.Item(27,4).Value = 219.1497099518948
$wsForecast.Cells.Item(28,1).Value = 45263.99999999999
$wsForecast.Cells.Item(28,2).Value = 95
$wsForecast.Cells.Item(28,3).Value = -23.99141889774057
$wsForecast.Cells.Item(28,4).Value = 206.5554522100581
$wsForecast.Cells.Item(29,1).Value = 45270.99999999999
$wsForecast.Cells.Item(29,2).Value = 95
$wsForecast.Cells.Item(29,3).Value = -23.356078549056
$wsForecast.Cells.Item(29,4).Value = 205.5351016526261
$wsForecast.Cells.Item(30,1).Value = 45277.99999999999
$wsForecast.Cells.Item(30,2).Value = 96
$wsForecast.Cells.Item(30,3).Value = -16.23737407454237
$wsForecast.Cells.Item(30,4).Value = 216.328801642275
$wsForecast.Cells.Item(31,1).Value = 45298.99999999999
$wsForecast.Cells.Item(31,2).Value = 98
$wsForecast.Cells.Item(31,3).Value = -16.63977432641435
$wsForecast.Cells.Item(31,4).Value = 215.6531347761043
$wsForecast.Cells.Item(32,1).Value = 45305.99999999999
$wsForecast.Cells.Item(32,2).Value = 99
$wsForecast.Cells.Item(32,3).Value = -19.97386671504709
$wsForecast.Cells.Item(32,4).Value = 216.2168644001829
$wsForecast.Cells.Item(33,1).Value = 45312.99999999999
$wsForecast.Cells.Item(33,2).Value = 100
$wsForecast.Cells.Item(33,3).Value = -18.86470009481863
$wsForecast.Cells.Item(33,4).Value = 226.8202506573602
$wsForecast.Cells.Item(34,1).Value = 45326.99999999999
$wsForecast.Cells.Item(34,2).Value = 101
$wsForecast.Cells.Item(34,3).Value = -21.82221912119683
$wsForecast.Cells.Item(34,4).Value = 219.2193967218491
$wsForecast.Cells.Item(35,1).Value = 45333.99999999999
$wsForecast.Cells.Item(35,2).Value = 102
$wsForecast.Cells.Item(35,3).Value = -24.92237487398979
$wsForecast.Cells.Item(35,4).Value = 229.5831000877806
$wsForecast.Cells.Item(36,1).Value = 45340.99999999999
$wsForecast.Cells.Item(36,2).Value = 102
$wsForecast.Cells.Item(36,3).Value = -14.92683461241772
$wsForecast.Cells.Item(36,4).Value = 221.9750159103768
$wsForecast.Cells.Item(37,1).Value = 45347.99999999999
$wsForecast.Cells.Item(37,2).Value = 103
$wsForecast.Cells.Item(37,3).Value = -18.96633905074599
$wsForecast.Cells.Item(37,4).Value = 225.0967669778581
$wsForecast.Cells.Item(38,1).Value = 45354.99999999999
$wsForecast.Cells.Item(38,2).Value = 104
$wsForecast.Cells.Item(38,3).Value = -12.13724036983847
$wsForecast.Cells.Item(38,4).Value = 222.0613202074793
$wsForecast.Cells.Item(39,1).Value = 45361.99999999999
$wsForecast.Cells.Item(39,2).Value = 105
$wsForecast.Cells.Item(39,3).Value = -2.391470068120325
$wsForecast.Cells.Item(39,4).Value = 228.8035896459875
$wsForecast.Cells.Item(40,1).Value = 45368.99999999999
$wsForecast.Cells.Item(40,2).Value = 105
$wsForecast.Cells.Item(40,3).Value = -11.94774539483964
$wsForecast.Cells.Item(40,4).Value = 225.3274109944802
$wsForecast.Cells.Item(41,1).Value = 45375.99999999999
$wsForecast.Cells.Item(41,2).Value = 106
$wsForecast.Cells.Item(41,3).Value = -7.242018474119892
$wsForecast.Cells.Item(41,4).Value = 222.5298076767683
$wsForecast.Cells.Item(42,1).Value = 45410.99999999999
$wsForecast.Cells.Item(42,2).Value = 110
$wsForecast.Cells.Item(42,3).Value = 1.568034056170704
$wsForecast.Cells.Item(42,4).Value = 228.7454947322329
$wsForecast.Cells.Item(43,1).Value = 45417.99999999999
$wsForecast.Cells.Item(43,2).Value = 110
$wsForecast.Cells.Item(43,3).Value = -14.2421104506997
$wsForecast.Cells.Item(43,4).Value = 218.286471791439
$wsForecast.Cells.Item(44,1).Value = 45424.99999999999
$wsForecast.Cells.Item(44,2).Value = 111
$wsForecast.Cells.Item(44,3).Value = -1.695307209440986
$wsForecast.Cells.Item(44,4).Value = 229.2473127085645
$wsForecast.Cells.Item(45,1).Value = 45431.99999999999
$wsForecast.Cells.Item(45,2).Value = 112
$wsForecast.Cells.Item(45,3).Value = -9.843584077693178
$wsForecast.Cells.Item(45,4).Value = 227.6438982320449
$wsForecast.Cells.Item(46,1).Value = 45438.99999999999
$wsForecast.Cells.Item(46,2).Value = 113
$wsForecast.Cells.Item(46,3).Value = -4.644662758345472
$wsForecast.Cells.Item(46,4).Value = 232.8107967009599
$wsForecast.Cells.Item(47,1).Value = 45445.99999999999
$wsForecast.Cells.Item(47,2).Value = 113
$wsForecast.Cells.Item(47,3).Value = -2.70268203555344
$wsForecast.Cells.Item(47,4).Value = 235.2945528143014
$wsForecast.Cells.Item(48,1).Value = 45452.99999999999
$wsForecast.Cells.Item(48,2).Value = 114
$wsForecast.Cells.Item(48,3).Value = -8.83894254665609
$wsForecast.Cells.Item(48,4).Value = 234.9429859684456
$wsForecast.Cells.Item(49,1).Value = 45459.99999999999
$wsForecast.Cells.Item(49,2).Value = 115
$wsForecast.Cells.Item(49,3).Value = -10.97545945619778
$wsForecast.Cells.Item(49,4).Value = 231.8052456055555
$wsForecast.Cells.Item(50,1).Value = 45466.99999999999
$wsForecast.Cells.Item(50,2).Value = 115
$wsForecast.Cells.Item(50,3).Value = -2.212963324059636
$wsForecast.Cells.Item(50,4).Value = 243.3435702694045
$wsForecast.Cells.Item(51,1).Value = 45473.99999999999
$wsForecast.Cells.Item(51,2).Value = 116
$wsForecast.Cells.Item(51,3).Value = 5.008082556072408
$wsForecast.Cells.Item(51,4).Value = 240.0596364721668
$wsForecast.Cells.Item(52,1).Value = 45487.99999999999
$wsForecast.Cells.Item(52,2).Value = 118
$wsForecast.Cells.Item(52,3).Value = -0.2367515049734939
$wsForecast.Cells.Item(52,4).Value = 238.4273431557592
$wsForecast.Cells.Item(53,1).Value = 45494.99999999999
$wsForecast.Cells.Item(53,2).Value = 118
$wsForecast.Cells.Item(53,3).Value = 2.888763509371894
$wsForecast.Cells.Item(53,4).Value = 231.9163434677177
$wsForecast.Cells.Item(54,1).Value = 45501.99999999999
$wsForecast.Cells.Item(54,2).Value = 119
$wsForecast.Cells.Item(54,3).Value = -2.886125558217324
$wsForecast.Cells.Item(54,4).Value = 238.9836546447158
$wsForecast.Cells.Item(55,1).Value = 45508.99999999999
$wsForecast.Cells.Item(55,2).Value = 120
$wsForecast.Cells.Item(55,3).Value = 4.883105621071559
$wsForecast.Cells.Item(55,4).Value = 239.3979266363435
$wsForecast.Cells.Item(56,1).Value = 45515.99999999999
$wsForecast.Cells.Item(56,2).Value = 120
$wsForecast.Cells.Item(56,3).Value = 5.169299296541038
$wsForecast.Cells.Item(56,4).Value = 226.7348062477021
$wsForecast.Cells.Item(57,1).Value = 45522.99999999999
$wsForecast.Cells.Item(57,2).Value = 121
$wsForecast.Cells.Item(57,3).Value = 3.4479898473424
$wsForecast.Cells.Item(57,4).Value = 236.8831724823323
$wsForecast.Cells.Item(58,1).Value = 45529.99999999999
$wsForecast.Cells.Item(58,2).Value = 122
$wsForecast.Cells.Item(58,3).Value = 6.477792511179454
$wsForecast.Cells.Item(58,4).Value = 239.9837712332994
$wsForecast.Cells.Item(59,1).Value = 45536.99999999999
$wsForecast.Cells.Item(59,2).Value = 123
$wsForecast.Cells.Item(59,3).Value = -2.043728295716493
$wsForecast.Cells.Item(59,4).Value = 240.3482646876423
$wsForecast.Cells.Item(60,1).Value = 45543.99999999999
$wsForecast.Cells.Item(60,2).Value = 123
$wsForecast.Cells.Item(60,3).Value = -0.7795899515261014
$wsForecast.Cells.Item(60,4).Value = 235.3351495830478
$wsForecast.Cells.Item(61,1).Value = 45550.99999999999
$wsForecast.Cells.Item(61,2).Value = 124
$wsForecast.Cells.Item(61,3).Value = 11.13971186748944
$wsForecast.Cells.Item(61,4).Value = 245.4952201491618
$wsForecast.Cells.Item(62,1).Value = 45557.99999999999
$wsForecast.Cells.Item(62,2).Value = 125
$wsForecast.Cells.Item(62,3).Value = 4.449100630257553
$wsForecast.Cells.Item(62,4).Value = 245.1091163468658
$wsForecast.Cells.Item(63,1).Value = 45564.99999999999
$wsForecast.Cells.Item(63,2).Value = 125
$wsForecast.Cells.Item(63,3).Value = 6.415670041365813
$wsForecast.Cells.Item(63,4).Value = 245.4078817415528
$wsForecast.Cells.Item(64,1).Value = 45571.99999999999
$wsForecast.Cells.Item(64,2).Value = 126
$wsForecast.Cells.Item(64,3).Value = 9.459672637105973
$wsForecast.Cells.Item(64,4).Value = 235.6895528487907
$wsForecast.Cells.Item(65,1).Value = 45578.99999999999
$wsForecast.Cells.Item(65,2).Value = 127
$wsForecast.Cells.Item(65,3).Value = 1.556720918510896
$wsForecast.Cells.Item(65,4).Value = 243.8776146193245
$wsForecast.Cells.Item(66,1).Value = 45585.99999999999
$wsForecast.Cells.Item(66,2).Value = 128
$wsForecast.Cells.Item(66,3).Value = 11.65126979727568
$wsForecast.Cells.Item(66,4).Value = 249.4428097737632
$wsForecast.Cells.Item(67,1).Value = 45592.99999999999
$wsForecast.Cells.Item(67,2).Value = 128
$wsForecast.Cells.Item(67,3).Value = 12.00076787072145
$wsForecast.Cells.Item(67,4).Value = 247.0902712030619
$wsForecast.Cells.Item(68,1).Value = 45599.99999999999
$wsForecast.Cells.Item(68,2).Value = 129
$wsForecast.Cells.Item(68,3).Value = 11.95488037289663
$wsForecast.Cells.Item(68,4).Value = 244.7349428280064
$wsForecast.Cells.Item(69,1).Value = 45606.99999999999
$wsForecast.Cells.Item(69,2).Value = 130
$wsForecast.Cells.Item(69,3).Value = 10.79645445491337
$wsForecast.Cells.Item(69,4).Value = 242.0534689861364
$wsForecast.Cells.Item(70,1).Value = 45613.99999999999
$wsForecast.Cells.Item(70,2).Value = 130
$wsForecast.Cells.Item(70,3).Value = 14.85274075996888
$wsForecast.Cells.Item(70,4).Value = 255.5552021949096
$wsForecast.Cells.Item(71,1).Value = 45620.99999999999
$wsForecast.Cells.Item(71,2).Value = 131
$wsForecast.Cells.Item(71,3).Value = 5.31786827367443
$wsForecast.Cells.Item(71,4).Value = 261.6445592451306
$wsForecast.Cells.Item(72,1).Value = 45627.99999999999
$wsForecast.Cells.Item(72,2).Value = 132
$wsForecast.Cells.Item(72,3).Value = 10.64822146721889
$wsForecast.Cells.Item(72,4).Value = 249.4224976383075
$wsForecast.Cells.Item(73,1).Value = 45634.99999999999
$wsForecast.Cells.Item(73,2).Value = 133
$wsForecast.Cells.Item(73,3).Value = 12.89664575967423
$wsForecast.Cells.Item(73,4).Value = 245.7573703852792
$wsForecast.Cells.Item(74,1).Value = 45641.99999999999
$wsForecast.Cells.Item(74,2).Value = 133
$wsForecast.Cells.Item(74,3).Value = 15.73385593680439
$wsForecast.Cells.Item(74,4).Value = 251.9560015880014
$wsForecast.Cells.Item(75,1).Value = 45648.99999999999
$wsForecast.Cells.Item(75,2).Value = 134
$wsForecast.Cells.Item(75,3).Value = 11.64715059569003
$wsForecast.Cells.Item(75,4).Value = 246.2752357065125
$wsForecast.Cells.Item(76,1).Value = 45655.99999999999
$wsForecast.Cells.Item(76,2).Value = 135
$wsForecast.Cells.Item(76,3).Value = 23.25475700524526
$wsForecast.Cells.Item(76,4).Value = 259.2296458517608

# --- 5. Formatting: reuse the existing header style (s=1) and date style (s=2) ---
# Copy header style (bold + border + centered) from an existing sheet header
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-time number-format style used for column A on existing sheets
$wsWeekly.Range("A2:A68").Copy()
$wsForecast.Range("A2:A76").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$wsForecast.Range("A1").Select()
